$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All price cells in column D are stored as text (inline strings) in the
# original workbook. Force NumberFormat to "@" (Text) before assignment so
# numeric-looking strings are not silently converted into numeric values,
# which preserves exact formatting (trailing zeros, precision, etc.).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.63"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "25.45"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.133"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05590"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.528"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.018"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8190"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8406"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1343"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06955"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02844"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001516"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005953"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006179"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.507"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03159"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.751"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04740"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1342"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001250"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004648"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009704"

$ws.Range("E27").Value = "26NitroExNTXBestin24h"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001389"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03662"

$ws.Range("B41").Value = "BKEXToken"

$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1364"

$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "KickToken"

$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006233"

$ws.Range("E42").Value = "41KickTokenKICK"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002641"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008300"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005302"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002117"
